$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23 - this shifts the existing rows 23:99 down
# to 24:100 and extends the used range to A1:R100, matching the diff.
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new record.
$ws.Cells.Item(23, 1).Value = 4
$ws.Cells.Item(23, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(23, 3).Value = "Los Lagos"
$ws.Cells.Item(23, 4).Value = 45071
$ws.Cells.Item(23, 5).Value = 10
$ws.Cells.Item(23, 6).Value = 100112043
$ws.Cells.Item(23, 7).Value = "Pepino dulce"
$ws.Cells.Item(23, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 30
$ws.Cells.Item(23, 11).Value = 21000
$ws.Cells.Item(23, 12).Value = 21000
$ws.Cells.Item(23, 13).Value = 21000
$ws.Cells.Item(23, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(23, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(23, 16).Value = 1167
$ws.Cells.Item(23, 17).Value = 18
$ws.Cells.Item(23, 18).Value = "Hortaliza"
